# Missouri_B team transition-probability matrix refresh.
# More simulated games + faster simulate-game logic changed the
# empirical transition frequencies baked into this sheet; only the
# probability values move, row totals still sum to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2021466905187835
$ws.Range("C2").Value = 0.5599284436493739
$ws.Range("J2").Value = 0.005366726296958855
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.07871198568872988

# Row 3
$ws.Range("B3").Value = 0.003058103975535168
$ws.Range("C3").Value = 0.02140672782874618
$ws.Range("J3").Value = 0.02446483180428135
$ws.Range("P3").Value = 0.7645259938837921
$ws.Range("S3").Value = 0.1865443425076453

# Row 4
$ws.Range("J4").Value = 0.05681818181818182
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.3068181818181818

# Row 5
$ws.Range("O5").Value = 0.2
$ws.Range("P5").Value = 0.8

# Row 6
$ws.Range("B6").Value = 0.06074766355140187
$ws.Range("D6").Value = 0.01635514018691589
$ws.Range("F6").Value = 0.06542056074766354
$ws.Range("J6").Value = 0.2383177570093458
$ws.Range("O6").Value = 0.01168224299065421
$ws.Range("Q6").Value = 0.1495327102803738
$ws.Range("R6").Value = 0.05841121495327103
$ws.Range("S6").Value = 0.3995327102803738

# Row 7
$ws.Range("B7").Value = 0.1064425770308123
$ws.Range("D7").Value = 0.01680672268907563
$ws.Range("F7").Value = 0.05042016806722689
$ws.Range("J7").Value = 0.1176470588235294
$ws.Range("O7").Value = 0.01400560224089636
$ws.Range("Q7").Value = 0.1764705882352941
$ws.Range("R7").Value = 0.09523809523809523
$ws.Range("S7").Value = 0.42296918767507

# Row 8
$ws.Range("B8").Value = 0.09773539928486293
$ws.Range("D8").Value = 0.02264600715137068
$ws.Range("E8").Value = 0.003575685339690107
$ws.Range("F8").Value = 0.05721096543504171
$ws.Range("J8").Value = 0.1156138259833135
$ws.Range("O8").Value = 0.01311084624553039
$ws.Range("Q8").Value = 0.1632896305125149
$ws.Range("R8").Value = 0.09535160905840286
$ws.Range("S8").Value = 0.4314660309892729

# Row 9
$ws.Range("B9").Value = 0.1027027027027027
$ws.Range("D9").Value = 0.01621621621621622
$ws.Range("F9").Value = 0.06216216216216217
$ws.Range("J9").Value = 0.1297297297297297
$ws.Range("O9").Value = 0.002702702702702703
$ws.Range("Q9").Value = 0.1540540540540541
$ws.Range("R9").Value = 0.0918918918918919
$ws.Range("S9").Value = 0.4405405405405405

# Row 10
$ws.Range("B10").Value = 0.1063917525773196
$ws.Range("D10").Value = 0.02268041237113402
$ws.Range("E10").Value = 0.0008247422680412372
$ws.Range("F10").Value = 0.07134020618556701
$ws.Range("J10").Value = 0.1331958762886598
$ws.Range("O10").Value = 0.01278350515463918
$ws.Range("Q10").Value = 0.2078350515463918
$ws.Range("R10").Value = 0.09072164948453608
$ws.Range("S10").Value = 0.3542268041237113

# Row 11
$ws.Range("G11").Value = 0.1295238095238095
$ws.Range("J11").Value = 0.09142857142857143
$ws.Range("K11").Value = 0.1828571428571429
$ws.Range("L11").Value = 0.5828571428571429
$ws.Range("S11").Value = 0.01333333333333333

# Row 12
$ws.Range("F12").Value = 0.003184713375796179
$ws.Range("G12").Value = 0.7834394904458599
$ws.Range("J12").Value = 0.1560509554140127
$ws.Range("K12").Value = 0.009554140127388535
$ws.Range("L12").Value = 0.01910828025477707
$ws.Range("S12").Value = 0.02866242038216561

# Row 13
$ws.Range("F13").Value = 0.01265822784810127
$ws.Range("G13").Value = 0.7341772151898734
$ws.Range("J13").Value = 0.2278481012658228
$ws.Range("S13").Value = 0.02531645569620253

# Row 14
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.0218978102189781
$ws.Range("H15").Value = 0.2214111922141119
$ws.Range("I15").Value = 0.05596107055961071
$ws.Range("J15").Value = 0.3600973236009732
$ws.Range("K15").Value = 0.06812652068126521
$ws.Range("M15").Value = 0.0072992700729927
$ws.Range("O15").Value = 0.08759124087591241
$ws.Range("S15").Value = 0.1776155717761557

# Row 16
$ws.Range("F16").Value = 0.02116402116402116
$ws.Range("H16").Value = 0.1931216931216931
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.4232804232804233
$ws.Range("K16").Value = 0.08994708994708994
$ws.Range("M16").Value = 0.01587301587301587
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.06613756613756613
$ws.Range("S16").Value = 0.1111111111111111

# Row 17
$ws.Range("F17").Value = 0.02078239608801956
$ws.Range("H17").Value = 0.176039119804401
$ws.Range("I17").Value = 0.09535452322738386
$ws.Range("J17").Value = 0.4205378973105134
$ws.Range("K17").Value = 0.08801955990220049
$ws.Range("M17").Value = 0.02322738386308068
$ws.Range("N17").Value = 0.001222493887530562
$ws.Range("O17").Value = 0.07946210268948656
$ws.Range("S17").Value = 0.09535452322738386

# Row 18
$ws.Range("F18").Value = 0.01269035532994924
$ws.Range("H18").Value = 0.1700507614213198
$ws.Range("I18").Value = 0.1116751269035533
$ws.Range("J18").Value = 0.3959390862944163
$ws.Range("K18").Value = 0.1040609137055838
$ws.Range("M18").Value = 0.01522842639593909
$ws.Range("O18").Value = 0.07868020304568528
$ws.Range("S18").Value = 0.1116751269035533

# Row 19
$ws.Range("F19").Value = 0.01806451612903226
$ws.Range("H19").Value = 0.2008602150537634
$ws.Range("I19").Value = 0.08516129032258064
$ws.Range("J19").Value = 0.3849462365591398
$ws.Range("K19").Value = 0.1053763440860215
$ws.Range("M19").Value = 0.02150537634408602
$ws.Range("N19").Value = 0.0004301075268817204
$ws.Range("O19").Value = 0.0675268817204301
$ws.Range("S19").Value = 0.1161290322580645
